$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 17
$ws.Range("H17").Value = 1311.1389
$ws.Range("J17").Value = 1314.3429
$ws.Range("L17").Value = 3943.0287
$ws.Range("N17").Value = -4279.028700000001

# row 21
$ws.Range("H21").Value = 45000
$ws.Range("I21").Value = 45000
$ws.Range("K21").Value = 45000
$ws.Range("M21").Value = -44532

# row 23
$ws.Range("H23").Value = 45000
$ws.Range("I23").Value = 45000
$ws.Range("K23").Value = 45000
$ws.Range("M23").Value = -44766

# row 29
$ws.Range("J29").Value = 999
$ws.Range("L29").Value = 2997
$ws.Range("N29").Value = -3559

# row 58
$ws.Range("H58").Value = 4028.4
$ws.Range("I58").Value = 1650
$ws.Range("J58").Value = 4893.273
$ws.Range("K58").Value = 4950
$ws.Range("L58").Value = 14679.819
$ws.Range("M58").Value = -4800
$ws.Range("N58").Value = -14979.819

# row 64
$ws.Range("H64").Value = 128000
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()

# row 67
$ws.Range("H67").Value = 128000
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()

# row 106
$ws.Range("H106").Value = 6330.3
$ws.Range("I106").Value = 7004
$ws.Range("K106").Value = 7004
$ws.Range("M106").Value = -6373

# row 132
$ws.Range("H132").Value = 4149
$ws.Range("I132").Value = 3951.8235
$ws.Range("K132").Value = 11855.4705
$ws.Range("M132").Value = -9325.470499999999

# row 141
$ws.Range("H141").Value = 7215.684
$ws.Range("I141").Value = 7439.933
$ws.Range("J141").Value = 6374.75
$ws.Range("K141").Value = 22319.799
$ws.Range("L141").Value = 19124.25
$ws.Range("M141").Value = -17139.799
$ws.Range("N141").Value = -29484.25

$ws = $wb.Worksheets.Item("ARM")
# row 5
$ws.Range("H5").Value = 4150
$ws.Range("I5").Value = 4150
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 4150
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -4038
$ws.Range("N5").ClearContents()

# row 32
$ws.Range("H32").Value = 7201
$ws.Range("I32").Value = 7143.8774
$ws.Range("K32").Value = 7143.8774
$ws.Range("M32").Value = -6856.8774

# row 113
$ws.Range("H113").Value = 124999.336
$ws.Range("J113").Value = 124999.336
$ws.Range("L113").Value = 124999.336
$ws.Range("N113").Value = -133677.336

# row 132
$ws.Range("H132").Value = 2352.6
$ws.Range("I132").Value = 2342.575
$ws.Range("J132").Value = 2432.8
$ws.Range("K132").Value = 7027.724999999999
$ws.Range("L132").Value = 7298.400000000001
$ws.Range("M132").Value = -4497.724999999999
$ws.Range("N132").Value = -12358.4

$ws = $wb.Worksheets.Item("BSM")
# row 4
$ws.Range("H4").Value = 4150
$ws.Range("I4").Value = 4150
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 4150
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -4035
$ws.Range("N4").ClearContents()

# row 86
$ws.Range("H86").Value = 5178.478
$ws.Range("I86").Value = 7348.5713
$ws.Range("J86").Value = 1802.7778
$ws.Range("K86").Value = 7348.5713
$ws.Range("L86").Value = 1802.7778
$ws.Range("M86").Value = -6225.5713
$ws.Range("N86").Value = -4048.7778

# row 89
$ws.Range("H89").Value = 5178.478
$ws.Range("I89").Value = 7348.5713
$ws.Range("J89").Value = 1802.7778
$ws.Range("K89").Value = 36742.85649999999
$ws.Range("L89").Value = 9013.889000000001
$ws.Range("M89").Value = -31126.85649999999
$ws.Range("N89").Value = -20245.889

# row 94
$ws.Range("H94").Value = 7622.1396
$ws.Range("J94").Value = 2426.4614
$ws.Range("L94").Value = 2426.4614
$ws.Range("N94").Value = -3328.4614

$ws = $wb.Worksheets.Item("CRP")
# row 31
$ws.Range("H31").Value = 8308.044
$ws.Range("I31").Value = 14005.7
$ws.Range("J31").Value = 3925.2307
$ws.Range("K31").Value = 14005.7
$ws.Range("L31").Value = 3925.2307
$ws.Range("M31").Value = -13710.7
$ws.Range("N31").Value = -4515.2307

# row 34
$ws.Range("H34").Value = 8308.044
$ws.Range("I34").Value = 14005.7
$ws.Range("J34").Value = 3925.2307
$ws.Range("K34").Value = 14005.7
$ws.Range("L34").Value = 3925.2307
$ws.Range("M34").Value = -13803.7
$ws.Range("N34").Value = -4329.2307

$ws = $wb.Worksheets.Item("CUL")
# row 4
$ws.Range("H4").Value = 42538120
$ws.Range("I4").Value = 44504252
$ws.Range("K4").Value = 133512756
$ws.Range("M4").Value = -133512644

# row 131
$ws.Range("H131").Value = 2723.8525
$ws.Range("J131").Value = 1968.341
$ws.Range("L131").Value = 5905.022999999999
$ws.Range("N131").Value = -15985.023

# row 132
$ws.Range("H132").Value = 22220.666
$ws.Range("I132").Value = 716.6
$ws.Range("K132").Value = 6449.400000000001
$ws.Range("M132").Value = -3919.400000000001

# row 133
$ws.Range("H133").Value = 14320.889
$ws.Range("I133").Value = 7224.5
$ws.Range("J133").Value = 19998
$ws.Range("K133").Value = 21673.5
$ws.Range("L133").Value = 59994
$ws.Range("M133").Value = -16613.5
$ws.Range("N133").Value = -70114

$ws = $wb.Worksheets.Item("GSM")
# row 103
$ws.Range("H103").Value = 29999
$ws.Range("J103").Value = 29999
$ws.Range("L103").Value = 29999
$ws.Range("N103").Value = -32343

# row 110
$ws.Range("H110").Value = 93555
$ws.Range("J110").Value = 93555
$ws.Range("L110").Value = 93555
$ws.Range("N110").Value = -101735

# row 126
$ws.Range("H126").Value = 6596.9355
$ws.Range("J126").Value = 3035.9
$ws.Range("L126").Value = 9107.700000000001
$ws.Range("N126").Value = -14047.7

# row 132
$ws.Range("H132").Value = 4358.175
$ws.Range("I132").Value = 4861
$ws.Range("K132").Value = 14583
$ws.Range("M132").Value = -12053

$ws = $wb.Worksheets.Item("LTW")
# row 55
$ws.Range("H55").Value = 4045
$ws.Range("I55").Value = 393.33334
$ws.Range("K55").Value = 393.33334
$ws.Range("M55").Value = -220.33334

# row 94
$ws.Range("H94").Value = 50000
$ws.Range("J94").Value = 50000
$ws.Range("L94").Value = 50000
$ws.Range("N94").Value = -51352

# row 100
$ws.Range("H100").Value = 21599.8
$ws.Range("I100").Value = 26333
$ws.Range("K100").Value = 26333
$ws.Range("M100").Value = -25792

$ws = $wb.Worksheets.Item("WVR")
# row 69
$ws.Range("H69").Value = 29999
$ws.Range("J69").Value = 29999
$ws.Range("L69").Value = 29999
$ws.Range("N69").Value = -31497

# row 72
$ws.Range("H72").Value = 29999
$ws.Range("J72").Value = 29999
$ws.Range("L72").Value = 89997
$ws.Range("N72").Value = -97485

# row 81
$ws.Range("H81").Value = 6541.7
$ws.Range("I81").Value = 7322.125
$ws.Range("J81").Value = 3420
$ws.Range("K81").Value = 14644.25
$ws.Range("L81").Value = 6840
$ws.Range("M81").Value = -13583.25
$ws.Range("N81").Value = -8962

# row 84
$ws.Range("H84").Value = 6541.7
$ws.Range("I84").Value = 7322.125
$ws.Range("J84").Value = 3420
$ws.Range("K84").Value = 73221.25
$ws.Range("L84").Value = 34200
$ws.Range("M84").Value = -67917.25
$ws.Range("N84").Value = -44808

# row 124
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
